$d = $word.ActiveDocument

# The resume's short form was missing the contact-info line that the
# long form has. Re-insert it as a new, centered paragraph directly
# below the name heading ("Dheeraj Chand"), matching the long-form
# layout.
#
# Using Find/Replace with a literal paragraph-mark ("^p") in the
# replacement text splits "Dheeraj Chand" into its own paragraph and
# creates a brand-new paragraph for the contact line. That new
# paragraph naturally inherits only the centered alignment (no bold /
# font-size run formatting), which matches the target markup exactly.
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
